$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Valor Mora" (F) and "Salario Basico" (G) values between
# row 16 (KEIDID MERCEDES TORRES ACOSTA) and row 19 (RAFAEL AUGUSTO BARRAZA RUIZ)
$f16 = $ws.Range("F16").Value2
$g16 = $ws.Range("G16").Value2
$f19 = $ws.Range("F19").Value2
$g19 = $ws.Range("G19").Value2

$ws.Range("F16").Value2 = $f19
$ws.Range("G16").Value2 = $g19
$ws.Range("F19").Value2 = $f16
$ws.Range("G19").Value2 = $g16
